$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.029.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.81%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.61'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.328'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.111.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.841.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.673'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.055.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('E22').Value = '  +2.83%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  +3.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.89'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.22%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.57%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.124'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0554'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  -1.42%  '
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('E34').Value = '  +23.10%  '
$ws.Range('E35').Value = '  +10.84%  '
$ws.Range('E36').Value = '  -4.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.752'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.19%  '
$ws.Range('E38').Value = '  +9.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '89.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0200'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.344.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.16%  '
$ws.Range('E45').Value = '  +3.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.84'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +79.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0529'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.026.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.41%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +15.36%  '
$ws.Range('E51').Value = '  +0.22%  '
